$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overall")

$ws.Range("H26").Value = 0.67239
$ws.Range("I26").Value = 0.019
$ws.Range("H27").Value = 0.23274
$ws.Range("I27").Value = 0.0432
$ws.Range("H28").Value = 0.6420400000000001
$ws.Range("I28").Value = 0.01542
$ws.Range("H29").Value = 0.08771
$ws.Range("I29").Value = 0.02748
$ws.Range("H30").Value = 0.67228
$ws.Range("I30").Value = 0.01826
$ws.Range("H31").Value = 0.23444
$ws.Range("I31").Value = 0.03833
$ws.Range("H32").Value = 0.64384
$ws.Range("I32").Value = 0.009220000000000001
$ws.Range("H33").Value = 0.09445000000000001
$ws.Range("I33").Value = 0.03298
$ws.Range("H34").Value = 0.68192
$ws.Range("I34").Value = 0.01761
$ws.Range("H35").Value = 0.15345
$ws.Range("I35").Value = 0.03596
$ws.Range("H36").Value = 0.68316
$ws.Range("I36").Value = 0.01759
$ws.Range("H37").Value = 0.172
$ws.Range("I37").Value = 0.0386
$ws.Range("H38").Value = 0.68192
$ws.Range("I38").Value = 0.01746
$ws.Range("H39").Value = 0.15345
$ws.Range("I39").Value = 0.03596
$ws.Range("H40").Value = 0.6836100000000001
$ws.Range("I40").Value = 0.01748
$ws.Range("H41").Value = 0.17032
$ws.Range("I41").Value = 0.03629
$ws.Range("H66").Value = 0.6561
$ws.Range("I66").Value = 0.01602
$ws.Range("H67").Value = 0.16871
$ws.Range("I67").Value = 0.03526
$ws.Range("H68").Value = 0.63214
$ws.Range("I68").Value = 0.01396
$ws.Range("H69").Value = 0.06244
$ws.Range("I69").Value = 0.03356
$ws.Range("H70").Value = 0.65565
$ws.Range("I70").Value = 0.01617
$ws.Range("H71").Value = 0.16871
$ws.Range("I71").Value = 0.03526
$ws.Range("H72").Value = 0.63282
$ws.Range("I72").Value = 0.01235
$ws.Range("H73").Value = 0.06748
$ws.Range("I73").Value = 0.02871
$ws.Range("H74").Value = 0.6761
$ws.Range("I74").Value = 0.01529
$ws.Range("H75").Value = 0.10798
$ws.Range("I75").Value = 0.02251
$ws.Range("H76").Value = 0.67667
$ws.Range("I76").Value = 0.01288
$ws.Range("H77").Value = 0.14675
$ws.Range("I77").Value = 0.01914
$ws.Range("H78").Value = 0.67588
$ws.Range("I78").Value = 0.01577
$ws.Range("H79").Value = 0.10798
$ws.Range("I79").Value = 0.02251
$ws.Range("H80").Value = 0.67655
$ws.Range("I80").Value = 0.01318
$ws.Range("H81").Value = 0.14843
$ws.Range("I81").Value = 0.01586
$ws.Range("H106").Value = 0.67687
$ws.Range("I106").Value = 0.01663
$ws.Range("H107").Value = 0.23613
$ws.Range("I107").Value = 0.03478
$ws.Range("H108").Value = 0.64911
$ws.Range("I108").Value = 0.01978
$ws.Range("H109").Value = 0.09107999999999999
$ws.Range("I109").Value = 0.02686
$ws.Range("H110").Value = 0.67653
$ws.Range("I110").Value = 0.01722
$ws.Range("H111").Value = 0.2378
$ws.Range("I111").Value = 0.0344
$ws.Range("H112").Value = 0.648
$ws.Range("I112").Value = 0.01956
$ws.Range("H113").Value = 0.09278
$ws.Range("I113").Value = 0.02721
$ws.Range("H114").Value = 0.69163
$ws.Range("I114").Value = 0.0131
$ws.Range("H115").Value = 0.15864
$ws.Range("I115").Value = 0.03531
$ws.Range("H116").Value = 0.69072
$ws.Range("I116").Value = 0.01498
$ws.Range("H117").Value = 0.18053
$ws.Range("I117").Value = 0.03975
$ws.Range("H118").Value = 0.69219
$ws.Range("I118").Value = 0.01357
$ws.Range("H119").Value = 0.15864
$ws.Range("I119").Value = 0.03531
$ws.Range("H120").Value = 0.69151
$ws.Range("I120").Value = 0.01621
$ws.Range("H121").Value = 0.18053
$ws.Range("I121").Value = 0.03975
$ws.Range("H146").Value = 0.66654
$ws.Range("I146").Value = 0.02356
$ws.Range("H147").Value = 0.19064
$ws.Range("I147").Value = 0.03682
$ws.Range("H148").Value = 0.63855
$ws.Range("I148").Value = 0.01922
$ws.Range("H149").Value = 0.06582
$ws.Range("I149").Value = 0.01817
$ws.Range("H150").Value = 0.66688
$ws.Range("I150").Value = 0.02318
$ws.Range("H151").Value = 0.18896
$ws.Range("I151").Value = 0.03608
$ws.Range("H152").Value = 0.64022
$ws.Range("I152").Value = 0.02175
$ws.Range("H153").Value = 0.08096
$ws.Range("I153").Value = 0.01474
$ws.Range("H154").Value = 0.68287
$ws.Range("I154").Value = 0.01804
$ws.Range("H155").Value = 0.11481
$ws.Range("I155").Value = 0.04334
$ws.Range("H156").Value = 0.68635
$ws.Range("I156").Value = 0.0186
$ws.Range("H157").Value = 0.16372
$ws.Range("I157").Value = 0.04149
$ws.Range("H158").Value = 0.68332
$ws.Range("I158").Value = 0.0188
$ws.Range("H159").Value = 0.11481
$ws.Range("I159").Value = 0.04334
$ws.Range("H160").Value = 0.6877
$ws.Range("I160").Value = 0.01761
$ws.Range("H161").Value = 0.16372
$ws.Range("I161").Value = 0.04183
